$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.527.63"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.619.28"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'521.38"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Value = "'144.49"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "2.627.23"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'6.26"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "3.078.35"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "58.573.26"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'20.67"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "2.621.64"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'344.68"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "'10.15"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'61.32"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'0.412"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "0.0₃0797"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "'6.18"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "'18.79"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'149.80"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'0.972"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "'36.58"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("E39").Value = "  -3.84%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'278.50"
$ws.Range("E42").Value = "  -4.32%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'0.0980"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "'19.55"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'0.0520"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "1.974.51"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "'4.62"
$ws.Range("E51").Value = "  -2.51%  "
